# Apply crypto price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain numeric-looking text (e.g. "27.986.21",
# "1.00", "0.0168") that must stay text, not be auto-coerced into numbers by
# Excel's usual text-to-number inference. Force text format before writing,
# then restore the default "Normal" style so no stray formatting is left behind.
$dRange = $ws.Range("D2","D3","D5","D8","D10","D12","D13","D16","D17","D18","D19","D20","D22","D23","D24","D25","D28","D29","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D46","D47","D49")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.986.21'
$ws.Range('E2').Value = '  +1.99%  '

$ws.Range('D3').Value = '1.642.06'
$ws.Range('E3').Value = '  +0.35%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = '212.89'
$ws.Range('E5').Value = '  +0.36%  '

$ws.Range('E6').Value = '  -1.10%  '

$ws.Range('E7').Value = '  -0.11%  '

$ws.Range('D8').Value = '23.31'
$ws.Range('E8').Value = '  +1.72%  '

$ws.Range('E9').Value = '  +2.68%  '

$ws.Range('D10').Value = '0.0612'
$ws.Range('E10').Value = '  +0.56%  '

$ws.Range('E11').Value = '  +0.75%  '

$ws.Range('D12').Value = '1.876.93'
$ws.Range('E12').Value = '  +0.44%  '

$ws.Range('D13').Value = '1.640.93'
$ws.Range('E13').Value = '  +0.21%  '

$ws.Range('E14').Value = '  +1.26%  '

$ws.Range('E15').Value = '  -3.71%  '

$ws.Range('D16').Value = '64.66'
$ws.Range('E16').Value = '  +0.80%  '

$ws.Range('D17').Value = '27.963.89'
$ws.Range('E17').Value = '  +2.01%  '

$ws.Range('D18').Value = '233.24'
$ws.Range('E18').Value = '  +1.81%  '

$ws.Range('D19').Value = '7.67'
$ws.Range('E19').Value = '  +2.08%  '

$ws.Range('D20').Value = '0.0₃0723'
$ws.Range('E20').Value = '  +0.21%  '

$ws.Range('E21').Value = '  -0.09%  '

$ws.Range('D22').Value = '4.32'
$ws.Range('E22').Value = '  +0.31%  '

$ws.Range('D23').Value = '9.98'
$ws.Range('E23').Value = '  +3.46%  '

$ws.Range('D24').Value = '2.08'
$ws.Range('E24').Value = '  +6.09%  '

$ws.Range('D25').Value = '150.40'
$ws.Range('E25').Value = '  +0.69%  '

$ws.Range('E26').Value = '  -0.35%  '

$ws.Range('E27').Value = '  -0.67%  '

$ws.Range('D28').Value = '15.70'
$ws.Range('E28').Value = '  +1.16%  '

$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.10%  '

$ws.Range('E30').Value = '  +0.34%  '

$ws.Range('E31').Value = '  -0.66%  '

$ws.Range('D32').Value = '3.31'
$ws.Range('E32').Value = '  +0.81%  '

$ws.Range('D33').Value = '1.471.39'
$ws.Range('E33').Value = '  +4.50%  '

$ws.Range('D34').Value = '3.11'
$ws.Range('E34').Value = '  -1.85%  '

$ws.Range('D35').Value = '1.55'
$ws.Range('E35').Value = '  -2.08%  '

$ws.Range('D36').Value = '2.36'
$ws.Range('E36').Value = '  -0.19%  '

$ws.Range('D37').Value = '0.567'
$ws.Range('E37').Value = '  -0.15%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '0.929'
$ws.Range('E38').Value = '  +13.23%  '

$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '0.881'
$ws.Range('E39').Value = '  +0.79%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.0168'
$ws.Range('E40').Value = '  +0.89%  '

$ws.Range('D41').Value = '69.63'
$ws.Range('E41').Value = '  +7.60%  '

$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.11%  '

$ws.Range('E43').Value = '  -1.98%  '

$ws.Range('E44').Value = '  -1.62%  '

$ws.Range('E45').Value = '  +0.39%  '

$ws.Range('D46').Value = '5.41'
$ws.Range('E46').Value = '  -1.31%  '

$ws.Range('D47').Value = '1.784.95'
$ws.Range('E47').Value = '  +0.41%  '

$ws.Range('E48').Value = '  +3.53%  '

$ws.Range('D49').Value = '86.40'
$ws.Range('E49').Value = '  +0.72%  '

$ws.Range('E50').Value = '  -0.13%  '

$ws.Range('E51').Value = '  +0.46%  '

# Restore default styling on the Price cells we just touched (keeps the
# worksheet's original "no explicit style" look while preserving the text).
$dRange.Style = "Normal"
